$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.248.21'
$ws.Range("D2").Style = $origStyle_D2
$ws.Range("E2").Value = '  +0.91%  '
$origStyle_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.304.72'
$ws.Range("D3").Style = $origStyle_D3
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  -0.01%  '
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.02'
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = '  +2.83%  '
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.88'
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = '  +1.82%  '
$origStyle_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("D7").Style = $origStyle_D7
$ws.Range("E7").Value = '  +1.93%  '
$ws.Range("E8").Value = '  -0.05%  '
$origStyle_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.304.72'
$ws.Range("D9").Style = $origStyle_D9
$ws.Range("E9").Value = '  +1.20%  '
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = '  +0.00%  '
$origStyle_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.85'
$ws.Range("D11").Style = $origStyle_D11
$ws.Range("E11").Value = '  +2.33%  '
$origStyle_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.402'
$ws.Range("D12").Style = $origStyle_D12
$ws.Range("E12").Value = '  +1.04%  '
$origStyle_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.881.23'
$ws.Range("D13").Style = $origStyle_D13
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("E14").Value = '  -2.36%  '
$origStyle_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.331.48'
$ws.Range("D15").Style = $origStyle_D15
$ws.Range("E15").Value = '  +0.87%  '
$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.65'
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E16").Value = '  +0.97%  '
$origStyle_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000163'
$ws.Range("D17").Style = $origStyle_D17
$ws.Range("E17").Value = '  +0.69%  '
$origStyle_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.295.40'
$ws.Range("D18").Style = $origStyle_D18
$ws.Range("E18").Value = '  +0.13%  '
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '424.91'
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = '  -2.32%  '
$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.49'
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = '  -1.05%  '
$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.08'
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = '  -0.33%  '
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.30'
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$origStyle_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = $origStyle_D23
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$origStyle_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.36'
$ws.Range("D24").Style = $origStyle_D24
$ws.Range("E24").Value = '  -1.39%  '
$origStyle_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.67'
$ws.Range("D25").Style = $origStyle_D25
$ws.Range("E25").Value = '  +0.11%  '
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.510'
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("E27").Value = '  +6.67%  '
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000114'
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = '  +2.05%  '
$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.44'
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = '  +7.01%  '
$origStyle_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = $origStyle_D30
$ws.Range("E30").Value = '  +0.13%  '
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.92'
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = '  -0.15%  '
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.32'
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("E33").Value = '  +0.05%  '
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.17'
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = '  +1.02%  '
$origStyle_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.59'
$ws.Range("D35").Style = $origStyle_D35
$ws.Range("E35").Value = '  +0.25%  '
$origStyle_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.19'
$ws.Range("D36").Style = $origStyle_D36
$ws.Range("E36").Value = '  +1.39%  '
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.59'
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = '  -0.06%  '
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.44'
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = '  +0.11%  '
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.872.20'
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = '  +3.62%  '
$origStyle_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.80'
$ws.Range("D40").Style = $origStyle_D40
$ws.Range("E40").Value = '  +1.44%  '
$origStyle_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.37'
$ws.Range("D41").Style = $origStyle_D41
$ws.Range("E41").Value = '  -0.92%  '
$origStyle_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.34'
$ws.Range("D42").Style = $origStyle_D42
$ws.Range("E42").Value = '  +0.80%  '
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.751'
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = '  -2.92%  '
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.69'
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = '  -1.23%  '
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.92'
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = '  -2.01%  '
$origStyle_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.31'
$ws.Range("D46").Style = $origStyle_D46
$ws.Range("E46").Value = '  +1.21%  '
$origStyle_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0640'
$ws.Range("D47").Style = $origStyle_D47
$ws.Range("E47").Value = '  -2.13%  '
$origStyle_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '315.02'
$ws.Range("D48").Style = $origStyle_D48
$ws.Range("E48").Value = '  -0.47%  '
$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.98'
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  +1.03%  '
$origStyle_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.102'
$ws.Range("D51").Style = $origStyle_D51
$ws.Range("E51").Value = '  +0.03%  '
